# tambah contoh di template
# Adds three example rows ("Contoh Ruangan 1/2/3") below the existing
# header row of the room template, widens column B to fit the new
# content, and moves the active selection to D3 (matching the saved
# state left behind by the author's Excel session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Contoh Ruangan 1"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Contoh Ruangan 2"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Contoh Ruangan 3"

# Widen column B so the longer example text fits (matches width="17").
$ws.Columns.Item(2).ColumnWidth = 16.14

# Leave the selection where the author's session left it.
[void]$ws.Range("D3").Select()
